$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2/B3/B4 hold numeric-looking "tracking order" IDs but must stay text
# (shared-string) cells with no style change, matching the original file.
# Typing the digits directly (Range.Value = "32384030") would be parsed as
# a number by Excel; going through a text *formula* + paste-values keeps
# the result a plain text cell without allocating a new quote-prefixed
# number-format style.
$scratch = $ws.Range("Z1")

$scratch.Formula = '="32384030"'
$scratch.Copy()
$ws.Range("B2").PasteSpecial(-4163)

$scratch.Formula = '="32384031"'
$scratch.Copy()
$ws.Range("B3").PasteSpecial(-4163)

$scratch.Formula = '="32384033"'
$scratch.Copy()
$ws.Range("B4").PasteSpecial(-4163)

$scratch.Clear()
$excel.CutCopyMode = $false
